$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'97.111.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "'3.706.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'2.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.54%  "
$ws.Range("D6").Value = "'235.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("D7").Value = "'656.39"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.431"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.78%  "
$ws.Range("D9").Value = "'1.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "'1.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").Value = "'3.705.56"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").Value = "'44.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "'0.0000309"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +14.99%  "
$ws.Range("D14").Value = "'0.206"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").Value = "'6.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "'4.404.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").Value = "'96.892.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "'9.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").Value = "'3.705.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").Value = "'13.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "'18.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("D22").Value = "'0.522"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").Value = "'524.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("D24").Value = "'3.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").Value = "'0.0000223"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.27%  "
$ws.Range("D26").Value = "'6.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.80%  "
$ws.Range("D27").Value = "'107.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.49%  "
$ws.Range("D28").Value = "'0.196"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +16.71%  "
$ws.Range("D29").Value = "'3.909.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").Value = "'13.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").Value = "'12.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").Value = "'3.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "'0.191"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.86%  "
$ws.Range("D35").Value = "'1.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.74%  "
$ws.Range("D36").Value = "'32.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("D38").Value = "'642.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.16%  "
$ws.Range("D39").Value = "'0.593"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").Value = "'8.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.05%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.166"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.499"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.40%  "
$ws.Range("D44").Value = "'6.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("D45").Value = "'2.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").Value = "'40.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.59%  "
$ws.Range("D47").Value = "'0.962"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "'0.0457"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").Value = "'2.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("D50").Value = "'23.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").Value = "'8.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.35%  "
